$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '62.619.33'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.437.75'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.14'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.19'
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = '2.434.23'
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.23'
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.78'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("D16").Value = '2.878.14'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '62.463.56'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").Value = '2.454.63'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.22'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.28'
$ws.Range("E20").Value = '  +3.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.59'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.06'
$ws.Range("E23").Value = '  +12.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.01'
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.03'
$ws.Range("E25").Value = '  -3.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '612.75'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.76'
$ws.Range("E27").Value = '  +1.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000101'
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("D29").Value = '2.558.93'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.49'
$ws.Range("E31").Value = '  +2.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.13'
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.139'
$ws.Range("E34").Value = '  -2.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.12'
$ws.Range("E35").Value = '  +4.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.51'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.377'
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.73'
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.33'
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '144.99'
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("E42").Value = '  -3.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.58'
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.86'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '147.62'
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.75'
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.86'
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0530'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.596'
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0231'
$ws.Range("E51").Value = '  -0.59%  '
